# adicionei a função flat field do Eduardo
# Recompute the R/G/B channel measurements on the "Color_smp" sheet
# after applying the flat-field correction.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Color_smp")

$ws.Range("A2").Value = 0.160625
$ws.Range("B2").Value = 0.1494485294117647
$ws.Range("C2").Value = 0.3810294117647058
$ws.Range("A3").Value = 0.2397009803921569
$ws.Range("B3").Value = 0.2279950980392157
$ws.Range("C3").Value = 0.4525612745098039
$ws.Range("A4").Value = 0.3130073529411764
$ws.Range("B4").Value = 0.3116225490196078
$ws.Range("C4").Value = 0.519892156862745
$ws.Range("A5").Value = 0.3857083333333333
$ws.Range("B5").Value = 0.3883897058823529
$ws.Range("C5").Value = 0.5952499999999999
$ws.Range("A6").Value = 0.4753504901960784
$ws.Range("B6").Value = 0.4857156862745098
$ws.Range("C6").Value = 0.6897009803921569
$ws.Range("A7").Value = 0.5735857843137255
$ws.Range("B7").Value = 0.5839509803921569
$ws.Range("C7").Value = 0.7953627450980392
$ws.Range("A8").Value = 0.2282352941176471
$ws.Range("B8").Value = 0.355125
$ws.Range("C8").Value = 0.337593137254902
$ws.Range("A9").Value = 0.2897794117647059
$ws.Range("B9").Value = 0.5289117647058824
$ws.Range("C9").Value = 0.7098848039215686
$ws.Range("A10").Value = 0.4487622549019608
$ws.Range("B10").Value = 0.2956421568627451
$ws.Range("C10").Value = 0.488953431372549
$ws.Range("A11").Value = 0.4752867647058824
$ws.Range("B11").Value = 0.374375
$ws.Range("C11").Value = 0.5948578431372549
$ws.Range("A12").Value = 0.3936495098039215
$ws.Range("B12").Value = 0.4515661764705883
$ws.Range("C12").Value = 0.4726936274509804
$ws.Range("A13").Value = 0.2216960784313725
$ws.Range("B13").Value = 0.3323700980392157
$ws.Range("C13").Value = 0.7474828431372549
$ws.Range("A14").Value = 0.2638578431372549
$ws.Range("B14").Value = 0.315093137254902
$ws.Range("C14").Value = 0.4995
$ws.Range("A15").Value = 0.4171691176470588
$ws.Range("B15").Value = 0.3069950980392157
$ws.Range("C15").Value = 0.7800637254901961
$ws.Range("A16").Value = 0.40875
$ws.Range("B16").Value = 0.2534607843137255
$ws.Range("C16").Value = 0.4465171568627451
$ws.Range("A17").Value = 0.4473970588235294
$ws.Range("B17").Value = 0.3688357843137255
$ws.Range("C17").Value = 0.5011838235294117
$ws.Range("A18").Value = 0.603406862745098
$ws.Range("B18").Value = 0.4089362745098039
$ws.Range("C18").Value = 0.3011446078431373
$ws.Range("A19").Value = 0.1967916666666667
$ws.Range("B19").Value = 0.5552598039215686
$ws.Range("C19").Value = 0.8262475490196078
$ws.Range("A20").Value = 0.3017794117647059
$ws.Range("B20").Value = 0.3170490196078431
$ws.Range("C20").Value = 0.5701299019607844
$ws.Range("A21").Value = 0.4174656862745098
$ws.Range("B21").Value = 0.4187279411764706
$ws.Range("C21").Value = 0.7986960784313725
$ws.Range("A22").Value = 0.361061274509804
$ws.Range("B22").Value = 0.289
$ws.Range("C22").Value = 0.5510980392156862
$ws.Range("A23").Value = 0.441468137254902
$ws.Range("B23").Value = 0.4523161764705882
$ws.Range("C23").Value = 0.7291544117647059
$ws.Range("A24").Value = 0.4766838235294117
$ws.Range("B24").Value = 0.4994289215686274
$ws.Range("C24").Value = 0.8183063725490196
$ws.Range("A25").Value = 0.2906470588235294
$ws.Range("B25").Value = 0.4766666666666667
$ws.Range("C25").Value = 0.8158455882352942
$ws.Range("A26").Value = 0.3066495098039216
$ws.Range("B26").Value = 0.3141372549019608
$ws.Range("C26").Value = 0.5956642156862745
$ws.Range("A27").Value = 0.3632450980392157
$ws.Range("B27").Value = 0.3806813725490196
$ws.Range("C27").Value = 0.8270196078431372
$ws.Range("A28").Value = 0.4167696078431373
$ws.Range("B28").Value = 0.4025294117647059
$ws.Range("C28").Value = 0.7010024509803922
$ws.Range("A29").Value = 0.4760906862745098
$ws.Range("B29").Value = 0.4688627450980392
$ws.Range("C29").Value = 0.7684607843137256
$ws.Range("A30").Value = 0.4597598039215686
$ws.Range("B30").Value = 0.5047377450980391
$ws.Range("C30").Value = 0.7937083333333333
$ws.Range("A31").Value = 0.438406862745098
$ws.Range("B31").Value = 0.5299338235294118
$ws.Range("C31").Value = 0.7946029411764706
